$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 633, shifting existing rows (633..679) down to (634..680)
$ws.Rows.Item(633).Insert()

# New row 633 duplicates the static columns of the (now-shifted) row 634,
# with its own Fecha/Volumen/Precio/Origen data (new weekly sample).
$ws.Cells.Item(633, 1).Value = 10
$ws.Cells.Item(633, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(633, 3).Value = "La Araucanía"
$ws.Cells.Item(633, 4).Value = 45265
$ws.Cells.Item(633, 5).Value = 9
$ws.Cells.Item(633, 6).Value = "Fruta"
$ws.Cells.Item(633, 7).Value = 100108
$ws.Cells.Item(633, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(633, 9).Value = 100108002
$ws.Cells.Item(633, 10).Value = "Mango"
$ws.Cells.Item(633, 11).Value = "Sin especificar"
$ws.Cells.Item(633, 12).Value = "Primera"
$ws.Cells.Item(633, 13).Value = 400
$ws.Cells.Item(633, 14).Value = 13000
$ws.Cells.Item(633, 15).Value = 13000
$ws.Cells.Item(633, 16).Value = 13000
$ws.Cells.Item(633, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(633, 18).Value = "Perú"
$ws.Cells.Item(633, 19).Value = 3250
$ws.Cells.Item(633, 20).Value = 4
